$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (as strings in column C to preserve text type)
$newRows = @(
    @(13, 1, "2024-06-15 04:15:51", 200, 6),
    @(14, 2, "2024-06-15 04:15:52", 200, 0)
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
